# "Added support for longer quotes, fixed surplus numnber"
#
# The surplus/markup multiplier in column K was set to 1.0565 on several
# quote-line rows; it should simply be 1 (no surplus) on those rows.
# Also move the active selection to C6 (where the user left off editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

# Fix the surplus number: K column multiplier 1.0565 -> 1 on the affected rows.
$surplusRows = @(16, 17, 20, 24, 27, 30, 31, 35)
foreach ($row in $surplusRows) {
    $ws.Cells.Item($row, 11).Value = 1
}

# Move the selection to C6, matching where the user is now working
# (part of the "support for longer quotes" editing pass).
$ws.Activate()
[void]$ws.Range("C6").Select()
